$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.076.95"
$ws.Range("E2").Value = "  +0.56%  "

# Row 3
$ws.Range("D3").Value = "2.543.05"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.580"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.91%  "

# Row 13
$ws.Range("E13").Value = "  +0.23%  "

# Row 14
$ws.Range("D14").Value = "2.932.92"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.10%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.514.40"
$ws.Range("E16").Value = "  -0.55%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "

# Row 18
$ws.Range("D18").Value = "43.050.85"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.22%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0972"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "

# Row 24
$ws.Range("E24").Value = "  +1.54%  "

# Row 25
$ws.Range("E25").Value = "  -0.35%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.02%  "

# Row 27
$ws.Range("E27").Value = "  +0.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.21%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.73%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0794"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.113"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.92%  "

# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("E41").Value = "  +10.13%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.78%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.28%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0305"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.90%  "

# Row 45
$ws.Range("E45").Value = "  +0.14%  "

# Row 46
$ws.Range("D46").Value = "2.029.96"
$ws.Range("E46").Value = "  -0.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "

# Row 49
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.30%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.785.65"
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
